# Insert a new vplan row just above the existing "Write to LFSR" row (row 30)
# to document a new "LFSR Bad Seed" test-plan item, per the xsecure vplan update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 30; everything below shifts down by one.
$ws.Rows.Item(30).Insert()

# The inserted row inherits default formatting; copy the "Feature" column
# formatting (borders/fill) from the row above so column B matches the
# surrounding block's look.
$ws.Range("B29").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's content (order matches how the strings were
# originally authored / appended to the shared string table).
$ws.Range("C30").Value = "LFSR Bad Seed"
$ws.Range("D30").Value = "Certain seeds should eventually lead to the LFSR locking up. (This is different from explicitly writing it to 0 and causing an immediate lockup.)"
$ws.Range("F30").Value = "Directed test"
$ws.Range("I30").Value = "TODO:WARNING:silabs-robin"
$ws.Range("E30").Value = "For all secureseed registers, write a value that eventually (but not immediately) leads to a lockup, and let the core run until that happens. (Could maybe be done with only a cover in formal, or with a simple directed test in sim, or both.)"
$ws.Range("G30").Value = "Directed Non-Self-Checking"
$ws.Range("H30").Value = "Functional coverage"

# Match the row height used for the new wrapped-text content.
$ws.Rows.Item(30).RowHeight = 97.2

# Leave the selection where the author left it after editing.
$ws.Range("E30").Select()
